$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.057.02'
$ws.Range("E2").Value = '  +4.45%  '
$ws.Range("D3").Value = '3.241.18'
$ws.Range("E3").Value = '  +2.24%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '577.07'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.96%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '176.04'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.25%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.604'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.65%  '
$ws.Range("D9").Value = '3.239.28'
$ws.Range("E9").Value = '  +2.13%  '
$ws.Range("E10").Value = '  +5.06%  '
$ws.Range("E11").Value = '  +1.43%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.407'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.44%  '
$ws.Range("D13").Value = '3.803.48'
$ws.Range("E13").Value = '  +2.20%  '
$ws.Range("E14").Value = '  +1.92%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.78'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.38%  '
$ws.Range("D16").Value = '67.029.27'
$ws.Range("E16").Value = '  +4.24%  '
$ws.Range("E17").Value = '  +4.14%  '
$ws.Range("D18").Value = '3.244.75'
$ws.Range("E18").Value = '  +2.41%  '
$ws.Range("E19").Value = '  +3.11%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.27'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.43%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '367.72'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.71%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.45'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.96%  '
$ws.Range("E23").Value = '  -0.08%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '70.10'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.71%  '
$ws.Range("B25").Value = 'WrappedeETH'
$ws.Range("C25").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D25").Value = '3.383.61'
$ws.Range("E25").Value = '  +2.12%  '
$ws.Range("B26").Value = 'Polygon'
$ws.Range("C26").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.505'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.13%  '
$ws.Range("E27").Value = '  +1.41%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.72'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.79%  '
$ws.Range("E29").Value = '  +1.44%  '
$ws.Range("E30").Value = '  +0.30%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.97'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.29%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.61'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.04%  '
$ws.Range("E33").Value = '  +1.93%  '
$ws.Range("B35").Value = 'Monero'
$ws.Range("C35").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '172.19'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +8.77%  '
$ws.Range("B36").Value = 'Fetch.AI'
$ws.Range("C36").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.23'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.87%  '
$ws.Range("B37").Value = 'Aptos'
$ws.Range("C37").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.76'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.96%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.50'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.43%  '
$ws.Range("E39").Value = '  +5.64%  '
$ws.Range("E40").Value = '  +11.14%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '26.67'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.63%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.59'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.13%  '
$ws.Range("D43").Value = '2.707.03'
$ws.Range("E43").Value = '  +2.55%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '6.36'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +6.35%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.28'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.88%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '40.33'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.48%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0670'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.09%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '24.51'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +5.42%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '333.16'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.79%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0279'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.98%  '
$ws.Range("E51").Value = '  +2.44%  '
